$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# MetDatas: add new MetStationName lookup values (Darwin / Alice
# Springs / Kalgoorlie) plus their MetStationId numbers. Do this
# before the ID->Id header renames below so the new strings are
# interned right after the still-referenced block.
# ------------------------------------------------------------------
$wsMet = $wb.Worksheets.Item("MetDatas")
$wsMet.Range("A2").Value = 1
$wsMet.Range("B2").Value = "Darwin"
$wsMet.Range("A3").Value = 2
$wsMet.Range("B3").Value = "Alice Springs"
$wsMet.Range("A4").Value = 3
$wsMet.Range("B4").Value = "Kalgoorlie"
$wsMet.Range("D3").Select() | Out-Null

# ------------------------------------------------------------------
# ExperimentInfos: header rename ExperimentInfoID -> ExperimentInfoId,
# plus sample rows.
# ------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("ExperimentInfos")
$wsExp.Range("A1").Value = "ExperimentInfoId"
$wsExp.Range("A2").Value = 1
$wsExp.Range("B2").Value = 1
$wsExp.Range("A3").Value = 2
$wsExp.Range("B3").Value = 1
$wsExp.Range("A4").Value = 3
$wsExp.Range("B4").Value = 2
$wsExp.Range("A5").Value = 4
$wsExp.Range("B5").Value = 2
$wsExp.Range("C12").Select() | Out-Null

# ------------------------------------------------------------------
# Stats: header rename StatsId -> StatId, plus sample rows.
# ------------------------------------------------------------------
$wsStats = $wb.Worksheets.Item("Stats")
$wsStats.Range("A1").Value = "StatId"
$wsStats.Range("A2").Value = 1
$wsStats.Range("B2").Value = 1
$wsStats.Range("C2").Value = 1
$wsStats.Range("D2").Value = 1
$wsStats.Range("A3").Value = 2
$wsStats.Range("B3").Value = 2
$wsStats.Range("C3").Value = 2
$wsStats.Range("D3").Value = 1
$wsStats.Range("A4").Value = 3
$wsStats.Range("B4").Value = 3
$wsStats.Range("C4").Value = 1
$wsStats.Range("D4").Value = 1
$wsStats.Range("A5").Value = 4
$wsStats.Range("B5").Value = 4
$wsStats.Range("C5").Value = 2
$wsStats.Range("D5").Value = 1
$wsStats.Range("B12").Select() | Out-Null

# ------------------------------------------------------------------
# SoilDatas: sample rows.
# ------------------------------------------------------------------
$wsSoil = $wb.Worksheets.Item("SoilDatas")
$wsSoil.Range("A2").Value = 1
$wsSoil.Range("B2").Value = 1
$wsSoil.Range("C2").Value = 1
$wsSoil.Range("A3").Value = 2
$wsSoil.Range("B3").Value = 1
$wsSoil.Range("C3").Value = 2
$wsSoil.Range("A4").Value = 3
$wsSoil.Range("B4").Value = 2
$wsSoil.Range("C4").Value = 1
$wsSoil.Range("A5").Value = 4
$wsSoil.Range("B5").Value = 2
$wsSoil.Range("C5").Value = 2
$wsSoil.Cells.Item(1, 1).Select() | Out-Null

# ------------------------------------------------------------------
# SoilLayerDatas: sample rows.
# ------------------------------------------------------------------
$wsLayer = $wb.Worksheets.Item("SoilLayerDatas")
$wsLayer.Range("A2").Value = 1
$wsLayer.Range("B2").Value = 1
$wsLayer.Range("C2").Value = 1
$wsLayer.Range("A3").Value = 2
$wsLayer.Range("B3").Value = 1
$wsLayer.Range("C3").Value = 2
$wsLayer.Range("A4").Value = 3
$wsLayer.Range("B4").Value = 2
$wsLayer.Range("C4").Value = 1
$wsLayer.Range("A5").Value = 4
$wsLayer.Range("B5").Value = 2
$wsLayer.Range("C5").Value = 2
$wsLayer.Range("E11").Select() | Out-Null

# Re-activate Stats as the visible sheet (matches the saved workbook's
# original active tab).
$wsStats.Activate() | Out-Null
